$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "all": add 2020/5/8 (serial 43959) row, bump prior-day total
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("all")

# Push the footer note (currently on row 31, col B only) down to row 32
$ws.Range("B31").Copy($ws.Range("B32"))
$ws.Range("B31").ClearContents()

# Yesterday's running test-count total ticks up by one
$ws.Range("B30").Value = 273

# New data row 31 - clone formatting from row 30, then overwrite values
$ws.Range("A30:H30").Copy($ws.Range("A31:H31"))
$ws.Range("A31").Value = 43959
$ws.Range("B31").Value = 273
$ws.Range("C31").Value = 268
$ws.Range("D31").Value = 79
$ws.Range("E31").Value = 69
$ws.Range("F31").Value = 10
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = 181

$ws.Activate()
[void]$ws.Range("A31").Select()

# -----------------------------------------------------------------
# Sheet "kobe": add matching 2020/5/8 row, fix yesterday's D/E values
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("kobe")

# Push the footer note (currently on row 86, col B only) down to row 87
$ws.Range("B86").Copy($ws.Range("B87"))
$ws.Range("B86").ClearContents()

# Correction to yesterday's figures
$ws.Range("D85").Value = 1
$ws.Range("E85").Value = 273

# New data row 86 - clone formatting from row 85, then overwrite values
$ws.Range("A85:J85").Copy($ws.Range("A86:J86"))
$ws.Range("A86").Value = 43959
$ws.Range("B86").Value = 0
$ws.Range("C86").Value = 2417
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 273
$ws.Range("F86").Value = 74
$ws.Range("G86").Value = 65
$ws.Range("H86").Value = 9
$ws.Range("I86").Value = 8
$ws.Range("J86").Value = 174

$ws.Activate()
[void]$ws.Range("A86").Select()

# -----------------------------------------------------------------
# Sheet "other": add matching 2020/5/8 row
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("other")

# Push the footer note (currently on row 61, col B only) down to row 62
$ws.Range("B61").Copy($ws.Range("B62"))
$ws.Range("B61").ClearContents()

# New data row 61 - clone formatting from row 60 (no column I there), then overwrite values
$ws.Range("A60:H60").Copy($ws.Range("A61:H61"))
$ws.Range("A61").Value = 43959
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = 12
$ws.Range("D61").Value = 5
$ws.Range("E61").Value = 4
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 7

$ws.Activate()
[void]$ws.Range("A61").Select()

# Restore "all" as the active/selected tab, matching the original workbook
$wb.Worksheets.Item("all").Activate()
